$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF holds a per-row "Date" value. Every data row (BF2:BF31) was
# populated with "5-31-2007-08" -- a mangled mm-dd-season string caused by
# how the NBA stats were scraped -- instead of the real game date. Correct
# it to the actual ISO date, "2008-05-31".
$dateRange = $ws.Range("BF2:BF31")

# Temporarily force Text format so Excel stores the literal string instead
# of re-parsing "2008-05-31" as a date serial number, then restore the
# default "Normal" cell style so no formatting residue is left behind on
# these cells (they were unstyled before the edit too).
$dateRange.NumberFormat = "@"
$dateRange.Value = "2008-05-31"
$dateRange.Style = "Normal"
